# Refresh Leve profit-calc columns (H:N) across the per-crafting-class
# sheets with latest Universalis average-price data. Generated by the
# scheduled price-update runner; values only (no formulas involved),
# columns: H currentAveragePrice, I currentAveragePriceNQ,
# J currentAveragePriceHQ, K LevePriceNQ, L LevePriceHQ,
# M LeveProfitNQ, N LeveProfitHQ.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 76.90000000000001  # row 42: Eye of the Beholder
$ws.Range("I42").Value = 47.375
$ws.Range("J42").Value = 195
$ws.Range("K42").Value = 142.125
$ws.Range("L42").Value = 585
$ws.Range("M42").Value = 87.875
$ws.Range("N42").Value = -1045

$ws.Range("H55").Value = 163.90909  # row 55: A Real Smooth Move
$ws.Range("I55").Value = 151.25
$ws.Range("J55").Value = 197.66667
$ws.Range("K55").Value = 151.25
$ws.Range("L55").Value = 197.66667
$ws.Range("M55").Value = 62.75
$ws.Range("N55").Value = -625.6666700000001

$ws.Range("H58").Value = 2020.8182  # row 58: A Matter of Vital Importance
$ws.Range("I58").Value = 1345.8
$ws.Range("J58").Value = 2583.3333
$ws.Range("K58").Value = 4037.4
$ws.Range("L58").Value = 7749.999899999999
$ws.Range("M58").Value = -3887.4
$ws.Range("N58").Value = -8049.999899999999

$ws.Range("H87").Value = 29333.334  # row 87: There Was a Late Fee
$ws.Range("J87").Value = 29333.334
$ws.Range("L87").Value = 29333.334
$ws.Range("N87").Value = -31829.334

$ws.Range("H90").Value = 29333.334  # row 90: A Gate Arcane Is Dragon's Bane (L)
$ws.Range("J90").Value = 29333.334
$ws.Range("L90").Value = 88000.00199999999
$ws.Range("N90").Value = -100480.002

$ws.Range("H129").Value = 740.6429000000001  # row 129: Practical Command
$ws.Range("I129").Value = 343.15
$ws.Range("J129").Value = 1102
$ws.Range("K129").Value = 1029.45
$ws.Range("L129").Value = 3306
$ws.Range("M129").Value = 3970.55
$ws.Range("N129").Value = -13306

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 119.14286  # row 4: Eyes Bigger than the Plate
$ws.Range("I4").Value = 119.14286
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 119.14286
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -3.142859999999999
$ws.Range("N4").ClearContents()

$ws.Range("H5").Value = 110  # row 5: The Alloyed Truth
$ws.Range("J5").Value = 100
$ws.Range("L5").Value = 100
$ws.Range("N5").Value = -324

$ws.Range("H17").Value = 0  # row 17: Cook Intentions
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H18").Value = 0  # row 18: Still the Best
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H63").Value = 2818.6365  # row 63: Rivets Run through It
$ws.Range("I63").Value = 2333.889
$ws.Range("K63").Value = 2333.889
$ws.Range("M63").Value = -1647.889

$ws.Range("H66").Value = 2818.6365  # row 66: A Riveting Revival (L)
$ws.Range("I66").Value = 2333.889
$ws.Range("K66").Value = 11669.445
$ws.Range("M66").Value = -8237.445

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 110  # row 4: Mending Fences
$ws.Range("J4").Value = 100
$ws.Range("L4").Value = 100
$ws.Range("N4").Value = -330

$ws.Range("H10").Value = 0  # row 10: Bring Me the Head Knife of Al'bedo Derssia
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()

$ws.Range("H17").Value = 0  # row 17: Peddle to the Metal
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H19").Value = 40000  # row 19: Twice as Slice
$ws.Range("J19").Value = 40000
$ws.Range("L19").Value = 40000
$ws.Range("N19").Value = -40346

$ws.Range("H55").Value = 29144.5  # row 55: Streamlining Operations
$ws.Range("J55").Value = 29144.5
$ws.Range("L55").Value = 29144.5
$ws.Range("N55").Value = -29690.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 0  # row 50: The Arsenal of Theocracy
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H107").Value = 1227.9231  # row 107: Built to Last
$ws.Range("I107").Value = 586.75
$ws.Range("J107").Value = 1777.5
$ws.Range("K107").Value = 586.75
$ws.Range("L107").Value = 1777.5
$ws.Range("M107").Value = 1333.25
$ws.Range("N107").Value = -5617.5

$ws.Range("H132").Value = 4120705  # row 132: Hull Lotta Damage
$ws.Range("I132").Value = 6862.4
$ws.Range("J132").Value = 9263008
$ws.Range("K132").Value = 20587.2
$ws.Range("L132").Value = 27789024
$ws.Range("M132").Value = -18057.2
$ws.Range("N132").Value = -27794084

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1725  # row 9: Jack of All Plates
$ws.Range("J9").Value = 1725
$ws.Range("L9").Value = 5175
$ws.Range("N9").Value = -5623

$ws.Range("H10").Value = 120.44444  # row 10: A Real Fungi
$ws.Range("I10").Value = 47
$ws.Range("J10").Value = 267.33334
$ws.Range("K10").Value = 141
$ws.Range("L10").Value = 802.0000200000001
$ws.Range("M10").Value = -2
$ws.Range("N10").Value = -1080.00002

$ws.Range("H12").Value = 264.96875  # row 12: Butter Me Up
$ws.Range("J12").Value = 217.16667
$ws.Range("L12").Value = 651.50001
$ws.Range("N12").Value = -997.50001

$ws.Range("H15").Value = 778.25  # row 15: Pretty Enough to Eat
$ws.Range("I15").Value = 375.66666
$ws.Range("J15").Value = 1019.8
$ws.Range("K15").Value = 1126.99998
$ws.Range("L15").Value = 3059.4
$ws.Range("M15").Value = -986.9999800000001
$ws.Range("N15").Value = -3339.4

$ws.Range("H16").Value = 1450.25  # row 16: Go Ahead and Dig In
$ws.Range("I16").Value = 925
$ws.Range("J16").Value = 1975.5
$ws.Range("K16").Value = 2775
$ws.Range("L16").Value = 5926.5
$ws.Range("M16").Value = -2602
$ws.Range("N16").Value = -6272.5

$ws.Range("H17").Value = 594.3333  # row 17: Chew the Fat
$ws.Range("I17").Value = 479.4
$ws.Range("K17").Value = 1438.2
$ws.Range("M17").Value = -1269.2

$ws.Range("H19").Value = 0  # row 19: The Bango Zango Diet
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()

$ws.Range("H20").Value = 3866.6667  # row 20: Omelette's Be Friends
$ws.Range("J20").Value = 5400
$ws.Range("L20").Value = 16200
$ws.Range("N20").Value = -16654

$ws.Range("H21").Value = 2750  # row 21: Shy Is the Oyster
$ws.Range("I21").Value = 3000
$ws.Range("J21").Value = 1000
$ws.Range("K21").Value = 9000
$ws.Range("L21").Value = 3000
$ws.Range("M21").Value = -8827
$ws.Range("N21").Value = -3346

$ws.Range("H22").Value = 1310  # row 22: A Total Nut Job
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1310
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3930
$ws.Range("N22").Value = -4268
$ws.Range("M22").ClearContents()

$ws.Range("H26").Value = 883.3333  # row 26: A Grape Idea
$ws.Range("I26").Value = 150
$ws.Range("K26").Value = 450
$ws.Range("M26").Value = -162

$ws.Range("H27").Value = 1310  # row 27: Brain Food
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1310
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 3930
$ws.Range("N27").Value = -4134
$ws.Range("M27").ClearContents()

$ws.Range("H32").Value = 3791.4285  # row 32: Convalescence Precedes Essence
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 3791.4285
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 11374.2855
$ws.Range("N32").Value = -11940.2855
$ws.Range("M32").ClearContents()

$ws.Range("H33").Value = 86.333336  # row 33: Cooking with Gas
$ws.Range("I33").Value = 90
$ws.Range("K33").Value = 540
$ws.Range("M33").Value = -257

$ws.Range("H54").Value = 5000  # row 54: Good Eats in Ishgard
$ws.Range("J54").Value = 5000
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -16118

$ws.Range("H69").Value = 1490  # row 69: Loving That Muffin Top
$ws.Range("I69").Value = 725
$ws.Range("J69").Value = 2000
$ws.Range("K69").Value = 2175
$ws.Range("L69").Value = 6000
$ws.Range("M69").Value = -1364
$ws.Range("N69").Value = -7622

$ws.Range("H72").Value = 1490  # row 72: Muffin of the Morn (L)
$ws.Range("I72").Value = 725
$ws.Range("J72").Value = 2000
$ws.Range("K72").Value = 6525
$ws.Range("L72").Value = 18000
$ws.Range("M72").Value = -2469
$ws.Range("N72").Value = -26112

$ws.Range("H131").Value = 859.0700000000001  # row 131: The Mountain Steeped
$ws.Range("I131").Value = 572.5
$ws.Range("J131").Value = 898.1477
$ws.Range("K131").Value = 1717.5
$ws.Range("L131").Value = 2694.4431
$ws.Range("M131").Value = 3322.5
$ws.Range("N131").Value = -12774.4431

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 91.07692  # row 2: Copper and Robbers
$ws.Range("I2").Value = 37.785713
$ws.Range("J2").Value = 153.25
$ws.Range("K2").Value = 37.785713
$ws.Range("L2").Value = 153.25
$ws.Range("M2").Value = 75.214287
$ws.Range("N2").Value = -379.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 5101.5  # row 32: Men Who Scare Up Goats
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 7169.1665
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 7169.1665
$ws.Range("M32").Value = -1683
$ws.Range("N32").Value = -7803.1665

$ws.Range("H33").Value = 0  # row 33: Just Rewards
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H34").Value = 12000  # row 34: Breeches Served Cold
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 12000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 12000
$ws.Range("N34").Value = -12344
$ws.Range("M34").ClearContents()

$ws.Range("H39").Value = 17866.666  # row 39: Quality over Quantity
$ws.Range("J39").Value = 17866.666
$ws.Range("L39").Value = 17866.666
$ws.Range("N39").Value = -18786.666

$ws.Range("H45").Value = 0  # row 45: Soft Shoe Shuffle
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H46").Value = 892.6667  # row 46: Supply Side Logic
$ws.Range("I46").Value = 833.3333
$ws.Range("J46").Value = 922.3333
$ws.Range("K46").Value = 833.3333
$ws.Range("L46").Value = 922.3333
$ws.Range("M46").Value = -645.3333
$ws.Range("N46").Value = -1298.3333

$ws.Range("H50").Value = 0  # row 50: The Birdmen of Ishgard
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H51").Value = 10999.5  # row 51: Skirt Chaser
$ws.Range("J51").Value = 10999.5
$ws.Range("L51").Value = 10999.5
$ws.Range("N51").Value = -11955.5

$ws.Range("H56").Value = 6250  # row 56: Hold On Tight
$ws.Range("I56").Value = 6250
$ws.Range("K56").Value = 6250
$ws.Range("M56").Value = -5559

$ws.Range("H57").Value = 14682  # row 57: Too Hot to Handle
$ws.Range("I57").Value = 8000
$ws.Range("J57").Value = 18023
$ws.Range("K57").Value = 8000
$ws.Range("L57").Value = 18023
$ws.Range("M57").Value = -7434
$ws.Range("N57").Value = -19155

$ws.Range("H58").Value = 26315  # row 58: Handle with Care
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 26315
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 26315
$ws.Range("N58").Value = -26835
$ws.Range("M58").ClearContents()
